# MODEL-INPUT CHANGES for vl, on second thought dont need this
#
# The "Constants" sheet had gained a row for
# "Number of VL tests recommended per person per year" (low/best/high =
# 1.5 / 2 / 2.5). On second thought that input isn't needed, so delete the
# entire row again; everything below it shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

$row = $ws.Rows.Item(41)
$row.Delete() | Out-Null

# Leave the selection on the row that now occupies the deleted row's
# position, matching what Excel does after an entire-row delete.
$ws.Rows.Item(41).Select() | Out-Null
